$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Paris 2023 Contenders Sticker Capsule"
$ws.Range("B2").Value = "$0.28 USD"

$ws.Range("A3").Value = "Paris 2023 Legends Sticker Capsule"
$ws.Range("B3").Value = "$0.28 USD"

$ws.Range("A4").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("B4").Value = "$0.28 USD"

$ws.Range("A5").Value = "Dreams & Nightmares Case"
$ws.Range("B5").Value = "$1.28 USD"

$ws.Range("A7").Value = "Paris 2023 Legends Sticker Capsule"
$ws.Range("B7").Value = "$0.28 USD"

$ws.Range("A8").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("B8").Value = "$0.28 USD"

$ws.Range("A9").Value = "Dreams & Nightmares Case"
$ws.Range("B9").Value = "$1.27 USD"

$ws.Range("A10").Value = "Paris 2023 Contenders Sticker Capsule"
$ws.Range("B10").Value = "$0.28 USD"

$ws.Range("A11").Value = "Paris 2023 Legends Sticker Capsule"
$ws.Range("B11").Value = "$0.28 USD"

$ws.Range("A12").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("B12").Value = "$0.28 USD"

$ws.Range("A13").Value = "Dreams & Nightmares Case"
$ws.Range("B13").Value = "$1.28 USD"

$wb.Save()
